# Hortaliza, Macroferia Regional de Talca - Berenjena
# A new weekly price record (week of 2021-08-30) was added to the
# dataset, inserted as a new row 35. All the following rows shift
# down by one (old row 35 -> new row 36, ..., old row 57 -> new row 58).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting existing rows 35-57 down to 36-58.
$ws.Rows("35:35").Insert()

# Populate the newly inserted row 35 with the new record's data.
$ws.Cells.Item(35, 1).Value = 5
$ws.Cells.Item(35, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(35, 3).Value = "Maule"
$ws.Cells.Item(35, 4).Value = 44438
$ws.Cells.Item(35, 5).Value = 7
$ws.Cells.Item(35, 6).Value = 100112001
$ws.Cells.Item(35, 7).Value = "Berenjena"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 300
$ws.Cells.Item(35, 11).Value = 6000
$ws.Cells.Item(35, 12).Value = 6000
$ws.Cells.Item(35, 13).Value = 6000
$ws.Cells.Item(35, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(35, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(35, 16).Value = 100
$ws.Cells.Item(35, 17).Value = 60
$ws.Cells.Item(35, 18).Value = "Hortaliza"
